$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.598.81'

$ws.Range('D3').Value = '1.958.93'
$ws.Range('E3').Value = '  +0.77%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.24'
$ws.Range('E5').Value = '  -0.09%  '

$ws.Range('E6').Value = '  +0.42%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.48'
$ws.Range('E7').Value = '  +4.81%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.375'
$ws.Range('E9').Value = '  +3.75%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0787'
$ws.Range('E10').Value = '  -7.01%  '

$ws.Range('E11').Value = '  +0.19%  '

$ws.Range('E12').Value = '  +5.26%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.248.08'
$ws.Range('E13').Value = '  +0.85%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.829'
$ws.Range('E14').Value = '  +1.90%  '

$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.68'
$ws.Range('E15').Value = '  +1.09%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.24'
$ws.Range('E16').Value = '  +1.22%  '

$ws.Range('D17').Value = '1.959.05'
$ws.Range('E17').Value = '  +0.70%  '

$ws.Range('D18').Value = '36.482.70'
$ws.Range('E18').Value = '  +0.20%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.58'
$ws.Range('E19').Value = '  +0.25%  '

$ws.Range('D20').Value = '0.0₃0850'
$ws.Range('E20').Value = '  -1.69%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '229.47'
$ws.Range('E21').Value = '  +0.25%  '

$ws.Range('E22').Value = '  +1.52%  '

$ws.Range('E23').Value = '  -0.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.43'
$ws.Range('E24').Value = '  +2.49%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  +1.76%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.142'
$ws.Range('E26').Value = '  +5.31%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.18'
$ws.Range('E27').Value = '  -0.49%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.38'
$ws.Range('E28').Value = '  -0.91%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.31'
$ws.Range('E29').Value = '  +0.48%  '

$ws.Range('E30').Value = '  +19.39%  '

$ws.Range('E31').Value = '  +0.79%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.75'
$ws.Range('E32').Value = '  +3.47%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0612'
$ws.Range('E33').Value = '  -0.98%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.47'
$ws.Range('E34').Value = '  +6.57%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.47'
$ws.Range('E35').Value = '  +9.06%  '

$ws.Range('E36').Value = '  -0.08%  '

$ws.Range('E37').Value = '  +3.94%  '

$ws.Range('E38').Value = '  -0.99%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.48'
$ws.Range('E39').Value = '  -12.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0963'
$ws.Range('E40').Value = '  -2.18%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.91'
$ws.Range('E41').Value = '  +0.68%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.16'
$ws.Range('E42').Value = '  +0.90%  '

$ws.Range('E43').Value = '  -0.07%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.81'
$ws.Range('E44').Value = '  -1.30%  '

$ws.Range('D45').Value = '1.362.05'
$ws.Range('E45').Value = '  +1.32%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.82'
$ws.Range('E46').Value = '  +2.59%  '

$ws.Range('E47').Value = '  -0.28%  '

$ws.Range('E48').Value = '  -1.30%  '

$ws.Range('E49').Value = '  +0.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.72'
$ws.Range('E50').Value = '  +5.22%  '

$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.55'
$ws.Range('E51').Value = '  +16.77%  '
